$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new expense row (row 28): 支出 生活费 400 on 2018-01-08, note "生活费(1/11-放寒假)" ---

# Copy the formatting from the row above (row 27) so the new row matches the
# existing striped/table look (fills, borders, number formats, alignment).
$ws.Range("B27:G27").Copy() | Out-Null
$ws.Range("B28:G28").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("B28").Value2 = 26
$ws.Range("C28").Value2 = "支出"
$ws.Range("D28").Value2 = 400
$ws.Range("E28").Value2 = 43108
$ws.Range("F28").Value2 = "生活费"
$ws.Range("G28").Value2 = "生活费(1/11-放寒假)"

# Update the active selection to reflect where the user ended up after the edit.
$ws.Activate()
$ws.Range("F22").Select() | Out-Null
